$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Delivery_results")

# New header cells for columns F and G
$ws.Range("F1").Value = "Average_waiting_time_(minutes)"
$ws.Range("G1").Value = "Average_queue_length"

# Copy style (border/font/alignment) from existing header cell E1 to F1:G1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing row values that changed
$ws.Range("D2").Value = 62
$ws.Range("E2").Value = 100

$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 76.67

# Fill in the new columns F and G for rows 2-6
$ws.Range("F2").Value = 0.63
$ws.Range("G2").Value = 6.08

$ws.Range("F3").Value = 0.23
$ws.Range("G3").Value = 2.04

$ws.Range("F4").Value = 0.2
$ws.Range("G4").Value = 1.74

$ws.Range("F5").Value = 0.19
$ws.Range("G5").Value = 1.72

$ws.Range("F6").Value = 0.29
$ws.Range("G6").Value = 2.43
